# Insert a new data row at row 178 on the (single) worksheet, pushing the
# existing rows 178-261 down to 179-262 (Excel's normal "insert row" shift
# behaviour), then populate the newly inserted row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the entire row shifts rows 178..261 down to 179..262 and grows
# the sheet's used range / dimension from R261 to R262 automatically.
$ws.Rows.Item(178).Insert()

# Populate the new row 178 with its data.
$ws.Range("A178").Value = 5
$ws.Range("B178").Value = "Macroferia Regional de Talca"
$ws.Range("C178").Value = "Maule"
$ws.Range("D178").Value = 44510
$ws.Range("E178").Value = 7
$ws.Range("F178").Value = 100112043
$ws.Range("G178").Value = "Pepino ensalada"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 500
$ws.Range("K178").Value = 8000
$ws.Range("L178").Value = 8000
$ws.Range("M178").Value = 8000
$ws.Range("N178").Value = "`$/caja 80 unidades"
$ws.Range("O178").Value = "Región del Maule"
$ws.Range("P178").Value = 100
$ws.Range("Q178").Value = 80
$ws.Range("R178").Value = "Hortaliza"
